$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-missing "url_pfadname" (English card name) values,
#     and fix a misspelled German card name. Order matches the original
#     authoring session so new shared-string entries line up the same way. ---
$ws.Range("C32").Value = 'Goblin Gang'
$ws.Range("B37").Value = 'Feuerwerkerin'
$ws.Range("C47").Value = 'Fisherman'
$ws.Range("C51").Value = 'Goblin Cage'
$ws.Range("C59").Value = 'Flying Machine'
$ws.Range("C72").Value = 'Goblin Demolisher'
$ws.Range("C79").Value = 'Golden Knight'
$ws.Range("C81").Value = 'Giant'
$ws.Range("C96").Value = 'Cannon Cart'
$ws.Range("C97").Value = 'Goblin Machine'
$ws.Range("C98").Value = 'Archer Queen'
$ws.Range("C99").Value = 'Goblinstein'
$ws.Range("C7").Value  = 'Goblins'
$ws.Range("C61").Value = 'Dark Prince'
$ws.Range("C95").Value = 'Executioner'
$ws.Range("C94").Value = 'Bowler'

# --- Correct "Anzahl" (count) for the Elixir Golem row ---
$ws.Range("E43").Value = 6

# --- Row 1 now uses the same explicit 15pt row height as the rest of the sheet ---
$ws.Rows.Item(1).RowHeight = 15

# --- Remove 202 unused trailing rows below the data (shifts the leftover
#     formatted-but-empty rows at the bottom of the sheet up accordingly) ---
$ws.Range("A102:A303").EntireRow.Delete()

# --- Reset scroll position / selection to the full data range ---
$ws.Range("A1:J101").Select()
